$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 26332.23
$ws.Range("I33").Value = 30562.395
$ws.Range("J33").Value = 3066.3333
$ws.Range("K33").Value = 30562.395
$ws.Range("L33").Value = 3066.3333
$ws.Range("M33").Value = -30333.395
$ws.Range("N33").Value = -3524.3333

$ws.Range("H131").Value = 3948.7307
$ws.Range("I131").Value = 480.55554
$ws.Range("J131").Value = 4674.628
$ws.Range("K131").Value = 1441.66662
$ws.Range("L131").Value = 14023.884
$ws.Range("M131").Value = 3598.33338
$ws.Range("N131").Value = -24103.884

$ws.Range("H132").Value = 4102129.5
$ws.Range("I132").Value = 4389971
$ws.Range("K132").Value = 13169913
$ws.Range("M132").Value = -13167383

$ws.Range("H137").Value = 1177.5897
$ws.Range("I137").Value = 739.7692
$ws.Range("J137").Value = 2053.2307
$ws.Range("K137").Value = 2219.3076
$ws.Range("L137").Value = 6159.6921
$ws.Range("M137").Value = 330.6923999999999
$ws.Range("N137").Value = -11259.6921

$ws.Range("H138").Value = 2082.1592
$ws.Range("I138").Value = 1471.2759
$ws.Range("J138").Value = 3263.2
$ws.Range("K138").Value = 4413.8277
$ws.Range("L138").Value = 9789.599999999999
$ws.Range("M138").Value = 726.1723000000002
$ws.Range("N138").Value = -20069.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19591.604
$ws.Range("I32").Value = 3570.5618
$ws.Range("K32").Value = 3570.5618
$ws.Range("M32").Value = -3283.5618

$ws.Range("H61").Value = 1237.2916
$ws.Range("I61").Value = 695.2857
$ws.Range("J61").Value = 2392
$ws.Range("K61").Value = 695.2857
$ws.Range("L61").Value = 2392
$ws.Range("M61").Value = -483.2857
$ws.Range("N61").Value = -2816

$ws.Range("H74").Value = 531.88
$ws.Range("I74").Value = 486.45456
$ws.Range("J74").Value = 865
$ws.Range("K74").Value = 486.45456
$ws.Range("L74").Value = 865
$ws.Range("M74").Value = 387.54544
$ws.Range("N74").Value = -2613

$ws.Range("H77").Value = 531.88
$ws.Range("I77").Value = 486.45456
$ws.Range("J77").Value = 865
$ws.Range("K77").Value = 2432.2728
$ws.Range("L77").Value = 4325
$ws.Range("M77").Value = 1935.7272
$ws.Range("N77").Value = -13061

$ws.Range("H97").Value = 38422.035
$ws.Range("I97").Value = 48583.332
$ws.Range("J97").Value = 2857.5
$ws.Range("K97").Value = 48583.332
$ws.Range("L97").Value = 2857.5
$ws.Range("M97").Value = -48087.332
$ws.Range("N97").Value = -3849.5

$ws.Range("H102").Value = 94169
$ws.Range("I102").Value = 170029.83
$ws.Range("J102").Value = 3136
$ws.Range("K102").Value = 170029.83
$ws.Range("L102").Value = 3136
$ws.Range("M102").Value = -168407.83
$ws.Range("N102").Value = -6380

$ws.Range("H132").Value = 1937.341
$ws.Range("I132").Value = 1980.481
$ws.Range("J132").Value = 1558.6666
$ws.Range("K132").Value = 5941.443
$ws.Range("L132").Value = 4675.9998
$ws.Range("M132").Value = -3411.443
$ws.Range("N132").Value = -9735.9998

$ws.Range("H136").Value = 1237.2916
$ws.Range("I136").Value = 695.2857
$ws.Range("J136").Value = 2392
$ws.Range("K136").Value = 2085.8571
$ws.Range("L136").Value = 7176
$ws.Range("M136").Value = 464.1428999999998
$ws.Range("N136").Value = -12276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 678.8889
$ws.Range("I94").Value = 397.5
$ws.Range("J94").Value = 904
$ws.Range("K94").Value = 397.5
$ws.Range("L94").Value = 904
$ws.Range("M94").Value = 53.5
$ws.Range("N94").Value = -1806

$ws.Range("H107").Value = 66725504
$ws.Range("I107").Value = 83405080
$ws.Range("J107").Value = 7193.6665
$ws.Range("K107").Value = 83405080
$ws.Range("L107").Value = 7193.6665
$ws.Range("M107").Value = -83403160
$ws.Range("N107").Value = -11033.6665

$ws.Range("H134").Value = 2541.8948
$ws.Range("I134").Value = 2205.1892
$ws.Range("K134").Value = 6615.567599999999
$ws.Range("M134").Value = -4080.567599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34718.773
$ws.Range("I31").Value = 1955.5333
$ws.Range("J31").Value = 51665.277
$ws.Range("K31").Value = 1955.5333
$ws.Range("L31").Value = 51665.277
$ws.Range("M31").Value = -1660.5333
$ws.Range("N31").Value = -52255.277

$ws.Range("H34").Value = 34718.773
$ws.Range("I34").Value = 1955.5333
$ws.Range("J34").Value = 51665.277
$ws.Range("K34").Value = 1955.5333
$ws.Range("L34").Value = 51665.277
$ws.Range("M34").Value = -1753.5333
$ws.Range("N34").Value = -52069.277

$ws.Range("H58").Value = 1171
$ws.Range("I58").Value = 1019.5476
$ws.Range("K58").Value = 1019.5476
$ws.Range("M58").Value = -816.5476

$ws.Range("H62").Value = 2662.5
$ws.Range("I62").Value = 2657.1428
$ws.Range("K62").Value = 2657.1428
$ws.Range("M62").Value = -2033.1428

$ws.Range("H65").Value = 2662.5
$ws.Range("I65").Value = 2657.1428
$ws.Range("K65").Value = 13285.714
$ws.Range("M65").Value = -10165.714

$ws.Range("H107").Value = 8138.7856
$ws.Range("J107").Value = 863
$ws.Range("L107").Value = 863
$ws.Range("N107").Value = -4703

$ws.Range("H122").Value = 860.375
$ws.Range("I122").Value = 1200
$ws.Range("J122").Value = 520.75
$ws.Range("K122").Value = 3600
$ws.Range("L122").Value = 1562.25
$ws.Range("M122").Value = -1150
$ws.Range("N122").Value = -6462.25

$ws.Range("H132").Value = 28304118
$ws.Range("I132").Value = 25002034
$ws.Range("K132").Value = 75006102
$ws.Range("M132").Value = -75003572

$ws.Range("H134").Value = 1132.5428
$ws.Range("I134").Value = 1063.2693
$ws.Range("J134").Value = 1332.6666
$ws.Range("K134").Value = 3189.8079
$ws.Range("L134").Value = 3997.9998
$ws.Range("M134").Value = -654.8078999999998
$ws.Range("N134").Value = -9067.9998

$ws.Range("H136").Value = 1171
$ws.Range("I136").Value = 1019.5476
$ws.Range("K136").Value = 3058.6428
$ws.Range("M136").Value = -508.6428000000001

$ws.Range("H140").Value = 60778.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 60778.5
$ws.Range("K140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -71138.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1175.7241
$ws.Range("I5").Value = 1145.2
$ws.Range("J5").Value = 1191.7894
$ws.Range("K5").Value = 3435.6
$ws.Range("L5").Value = 3575.3682
$ws.Range("M5").Value = -3323.6
$ws.Range("N5").Value = -3799.3682

$ws.Range("H131").Value = 8483.616
$ws.Range("J131").Value = 8766.786
$ws.Range("L131").Value = 26300.358
$ws.Range("N131").Value = -36380.358

$ws.Range("H135").Value = 1175.7241
$ws.Range("I135").Value = 1145.2
$ws.Range("J135").Value = 1191.7894
$ws.Range("K135").Value = 10306.8
$ws.Range("L135").Value = 10726.1046
$ws.Range("M135").Value = -7771.800000000001
$ws.Range("N135").Value = -15796.1046

$ws.Range("H139").Value = 1818.421
$ws.Range("I139").Value = 1018.5238
$ws.Range("J139").Value = 2806.5293
$ws.Range("K139").Value = 3055.5714
$ws.Range("L139").Value = 8419.5879
$ws.Range("M139").Value = 2084.4286
$ws.Range("N139").Value = -18699.5879

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 91004350
$ws.Range("I80").Value = 166839650
$ws.Range("J80").Value = 1999.8
$ws.Range("K80").Value = 166839650
$ws.Range("L80").Value = 1999.8
$ws.Range("M80").Value = -166838652
$ws.Range("N80").Value = -3995.8

$ws.Range("H83").Value = 91004350
$ws.Range("I83").Value = 166839650
$ws.Range("J83").Value = 1999.8
$ws.Range("K83").Value = 834198250
$ws.Range("L83").Value = 9999
$ws.Range("M83").Value = -834193258
$ws.Range("N83").Value = -19983

$ws.Range("H102").Value = 252265.8
$ws.Range("I102").Value = 1643.4166
$ws.Range("J102").Value = 502888.16
$ws.Range("K102").Value = 1643.4166
$ws.Range("L102").Value = 502888.16
$ws.Range("M102").Value = -21.41660000000002
$ws.Range("N102").Value = -506132.16

$ws.Range("H107").Value = 1122995
$ws.Range("I107").Value = 732.5
$ws.Range("J107").Value = 2020805
$ws.Range("K107").Value = 732.5
$ws.Range("L107").Value = 2020805
$ws.Range("M107").Value = 1187.5
$ws.Range("N107").Value = -2024645

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2117.7234
$ws.Range("I132").Value = 2249.0232
$ws.Range("J132").Value = 706.25
$ws.Range("K132").Value = 6747.069600000001
$ws.Range("L132").Value = 2118.75
$ws.Range("M132").Value = -4217.069600000001
$ws.Range("N132").Value = -7178.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 36000
$ws.Range("J57").Value = 36000
$ws.Range("L57").Value = 36000
$ws.Range("N57").Value = -37508

$ws.Range("H113").Value = 713.95654
$ws.Range("I113").Value = 567.2143
$ws.Range("J113").Value = 942.2222
$ws.Range("K113").Value = 1701.6429
$ws.Range("L113").Value = 2826.6666
$ws.Range("M113").Value = 468.3571000000002
$ws.Range("N113").Value = -7166.6666

$ws.Range("H122").Value = 1975.3928
$ws.Range("I122").Value = 1565.6471
$ws.Range("K122").Value = 4696.9413
$ws.Range("M122").Value = -2246.9413

$ws.Range("H136").Value = 692.11365
$ws.Range("I136").Value = 430.26666
$ws.Range("K136").Value = 1290.79998
$ws.Range("M136").Value = 1259.20002
